$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("imgs")

# New "Purple Swirl" source row (row 6): Content | Source | Creator | URL
$ws.Range("A6").Value = "Purple Swirl"
$ws.Range("B6").Value = "Unsplash"
$ws.Range("C6").Value = "Pawel Czerwinski"
$ws.Range("D6").Value = "https://unsplash.com/photos/fPN1w7bIuNU"

# Hyperlink the Creator cell to the photographer's Unsplash profile,
# matching the existing "Martin Katler" (C4) hyperlink pattern.
$ws.Hyperlinks.Add($ws.Range("C6"), "https://unsplash.com/@pawel_czerwinski", "", "https://unsplash.com/@pawel_czerwinski", "Pawel Czerwinski")

# Restore the cell style that Hyperlinks.Add overwrote, copying the same
# "Creator+hyperlink" look already used by C4.
$ws.Range("C4").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Final selection state left by the editing session.
$ws.Range("A1:XFD1048576").Select()
